# Updated transition-probability matrix values on Sheet1 (NJIT_B team matrix)
# following refreshed game simulation data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2195121951219512
$ws.Range("C2").Value = 0.4908536585365854
$ws.Range("J2").Value = 0.03353658536585366
$ws.Range("P2").Value = 0.1432926829268293
$ws.Range("S2").Value = 0.1128048780487805
$ws.Range("B3").Value = 0.01807228915662651
$ws.Range("C3").Value = 0.01204819277108434
$ws.Range("J3").Value = 0.02409638554216868
$ws.Range("P3").Value = 0.7168674698795181
$ws.Range("S3").Value = 0.2289156626506024
$ws.Range("J4").Value = 0.06
$ws.Range("P4").Value = 0.7
$ws.Range("S4").Value = 0.24
$ws.Range("B6").Value = 0.05633802816901409
$ws.Range("D6").Value = 0.004694835680751174
$ws.Range("F6").Value = 0.06103286384976526
$ws.Range("J6").Value = 0.3051643192488263
$ws.Range("O6").Value = 0.0187793427230047
$ws.Range("Q6").Value = 0.2065727699530517
$ws.Range("R6").Value = 0.06572769953051644
$ws.Range("S6").Value = 0.2816901408450704
$ws.Range("B7").Value = 0.1494252873563219
$ws.Range("D7").Value = 0.02873563218390805
$ws.Range("F7").Value = 0.04597701149425287
$ws.Range("J7").Value = 0.1436781609195402
$ws.Range("O7").Value = 0.005747126436781609
$ws.Range("Q7").Value = 0.1781609195402299
$ws.Range("R7").Value = 0.09195402298850575
$ws.Range("S7").Value = 0.3563218390804598
$ws.Range("B8").Value = 0.1293800539083558
$ws.Range("D8").Value = 0.02425876010781671
$ws.Range("F8").Value = 0.06199460916442048
$ws.Range("J8").Value = 0.1320754716981132
$ws.Range("O8").Value = 0.008086253369272238
$ws.Range("Q8").Value = 0.1455525606469003
$ws.Range("R8").Value = 0.1078167115902965
$ws.Range("S8").Value = 0.3908355795148248
$ws.Range("B9").Value = 0.07909604519774012
$ws.Range("D9").Value = 0.01129943502824859
$ws.Range("F9").Value = 0.05084745762711865
$ws.Range("J9").Value = 0.1355932203389831
$ws.Range("O9").Value = 0.005649717514124294
$ws.Range("Q9").Value = 0.2542372881355932
$ws.Range("R9").Value = 0.05649717514124294
$ws.Range("S9").Value = 0.4067796610169492
$ws.Range("B10").Value = 0.1248959200666112
$ws.Range("D10").Value = 0.02830974188176519
$ws.Range("E10").Value = 0.0008326394671107411
$ws.Range("F10").Value = 0.07327227310574522
$ws.Range("J10").Value = 0.1232306411323897
$ws.Range("O10").Value = 0.009159034138218152
$ws.Range("Q10").Value = 0.2014987510407993
$ws.Range("R10").Value = 0.07660283097418817
$ws.Range("S10").Value = 0.3621981681931724
$ws.Range("G11").Value = 0.1473684210526316
$ws.Range("J11").Value = 0.0912280701754386
$ws.Range("K11").Value = 0.2245614035087719
$ws.Range("L11").Value = 0.5087719298245614
$ws.Range("S11").Value = 0.02807017543859649
$ws.Range("G12").Value = 0.7450980392156863
$ws.Range("J12").Value = 0.196078431372549
$ws.Range("K12").Value = 0.006535947712418301
$ws.Range("L12").Value = 0.0392156862745098
$ws.Range("S12").Value = 0.0130718954248366
$ws.Range("G13").Value = 0.7142857142857143
$ws.Range("J13").Value = 0.1785714285714286
$ws.Range("S13").Value = 0.1071428571428571
$ws.Range("F15").Value = 0.04736842105263158
$ws.Range("H15").Value = 0.1052631578947368
$ws.Range("I15").Value = 0.08947368421052632
$ws.Range("J15").Value = 0.4157894736842105
$ws.Range("K15").Value = 0.07894736842105263
$ws.Range("M15").Value = 0.02105263157894737
$ws.Range("O15").Value = 0.06842105263157895
$ws.Range("S15").Value = 0.1736842105263158
$ws.Range("F16").Value = 0.01047120418848168
$ws.Range("H16").Value = 0.193717277486911
$ws.Range("I16").Value = 0.06806282722513089
$ws.Range("J16").Value = 0.418848167539267
$ws.Range("K16").Value = 0.1151832460732984
$ws.Range("M16").Value = 0.02094240837696335
$ws.Range("N16").Value = 0.005235602094240838
$ws.Range("O16").Value = 0.02617801047120419
$ws.Range("S16").Value = 0.1413612565445026
$ws.Range("F17").Value = 0.02450980392156863
$ws.Range("H17").Value = 0.1691176470588235
$ws.Range("I17").Value = 0.08333333333333333
$ws.Range("J17").Value = 0.4411764705882353
$ws.Range("K17").Value = 0.1053921568627451
$ws.Range("M17").Value = 0.01225490196078431
$ws.Range("N17").Value = 0.002450980392156863
$ws.Range("O17").Value = 0.05147058823529412
$ws.Range("S17").Value = 0.1102941176470588
$ws.Range("F18").Value = 0.01176470588235294
$ws.Range("H18").Value = 0.1588235294117647
$ws.Range("I18").Value = 0.07647058823529412
$ws.Range("J18").Value = 0.4235294117647059
$ws.Range("K18").Value = 0.07058823529411765
$ws.Range("O18").Value = 0.1058823529411765
$ws.Range("S18").Value = 0.1529411764705882
$ws.Range("F19").Value = 0.02264808362369338
$ws.Range("H19").Value = 0.1898954703832753
$ws.Range("I19").Value = 0.08797909407665505
$ws.Range("J19").Value = 0.3632404181184669
$ws.Range("K19").Value = 0.1088850174216028
$ws.Range("M19").Value = 0.01393728222996516
$ws.Range("N19").Value = 0.0008710801393728223
$ws.Range("O19").Value = 0.07926829268292683
$ws.Range("S19").Value = 0.1332752613240418
